$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns retain their original text formatting
# so numeric-looking strings (e.g. "0.110", "42.714.00") are not coerced into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.714.00"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.353.50"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.31"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.85"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.21%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.43"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.22%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.14%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.17"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.709.06"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.350.39"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.607.06"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.89"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.82%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.54"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.67"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "263.47"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.04"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.19%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.46"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.75"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.19"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.80%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.25"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.05"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.34%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("B36").Value = "Kaspa"

$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.17%  "

$ws.Range("B37").Value = "RenderToken"

$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.20%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.78"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.42%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.237"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.85"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.56"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.05%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.81"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +22.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.86"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.90%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.14"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.57%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.26%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.04%  "
